# Comment out unnecessary log messages for Append .NET/JAVA .xaml
#
# The Trello backlog dispatcher re-ran and picked up fresh board ids for the
# Java/.NET foundations boards. Refresh the board ids captured in column A,
# and append the row for Moiya Josephs (previously missing its own board id
# - it had been sharing the Marielle Nolasco board id/name by mistake).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Board id refreshed for Andrew Shields-Java Foundations Project.
$ws.Range("A3").Value = "62b8ebddeb5f78389b68cbbf"

# Board id refreshed for Marielle Nolasco-.NET Foundations Project.
$ws.Range("A4").Value = "62b8d257b8f7598718367a01"

# Row 5 picked up a new board id this run, still logged against the
# Marielle Nolasco-.NET Foundations Project board name.
$ws.Range("A5").Value = "62b8ebddcf728a4dbd624d11"
$ws.Range("B5").Value = "Marielle Nolasco-.NET Foundations Project"

# Append the corrected entry for Moiya Josephs with her own board id.
$ws.Range("A6").Value = "62b8ebdd5226315c99f9b256"
$ws.Range("B6").Value = "Moiya Josephs-Java Foundations Project"

# Match the formatting of the rest of the table for the new row.
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)

# Sheet's default column width widened to fit the longer board ids/names.
$ws.StandardWidth = 167.872656
